$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 15.01856033333333
$ws.Range("N2").Value = 45.055681
$ws.Range("O2").Value = 0.4908713633047416
$ws.Range("P2").Value = 0.4908713633047417
$ws.Range("Q2").Value = 15.47937982622778
$ws.Range("R2").Value = 139.31441843605
$ws.Range("S2").Value = 0.4908713633047416
$ws.Range("T2").Value = 0.4908713633047417

# Row 3 updates
$ws.Range("O3").Value = 0.3099803572711625
$ws.Range("P3").Value = 0.3099803572711625
$ws.Range("Q3").Value = 9.775073568288889
$ws.Range("R3").Value = 87.9756621146
$ws.Range("S3").Value = 0.3099803572711625
$ws.Range("T3").Value = 0.3099803572711625

# Row 4 updates
$ws.Range("O4").Value = 0.1991482794240958
$ws.Range("P4").Value = 0.1991482794240958
$ws.Range("R4").Value = 56.5203611466
$ws.Range("S4").Value = 0.1991482794240958
$ws.Range("T4").Value = 0.1991482794240958
